$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new survey entry ("明略科技") was added as the new first data row (row 2).
# This pushes every existing data row (2-25) down by one (to 3-26).
# Column A holds a simple positional index (0,1,2,...) per row, so after the
# shift it keeps counting up naturally; we just need to add the new A26 = 24
# for the row that now lands at the bottom.

# 1) Shift all existing data (columns B:S, rows 2-25) down by one row.
$srcVals = $ws.Range("B2:S25").Value()

# A handful of source cells hold plain text that LOOKS like a number/percent
# ("995", "100%", "8%"). Re-assigning a bare .Value would let Excel's normal
# type-inference silently convert those into numeric/percentage cells, which
# would change their stored type. Prefix them with a literal leading
# apostrophe (exactly what typing '995 into Excel does) so they round-trip
# as plain text, matching the source.
$srcVals[20,7]  = "'" + $srcVals[20,7]    # H21 -> "995"
$srcVals[7,18]  = "'" + $srcVals[7,18]    # S8  -> "100%"
$srcVals[16,17] = "'" + $srcVals[16,17]   # R17 -> "8%"

$ws.Range("B3:S26").Value = $srcVals

# 2) Give the newly created row 26 the same cell style as the row above it
#    (column A index cells use style "1": bold/centered/bordered), then set
#    its numeric index value.
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A26").Value = 24

# 3) Fill in row 2 with the brand-new company's survey data.
$ws.Range("B2").Value = "明略科技"
$ws.Range("C2").Value = "云密城J"
$ws.Range("D2").Value = "营销智能"
$ws.Range("E2").Value = "Java"
$ws.Range("F2").Value = "9:30-18:30"
$ws.Range("G2").Value = "1.5h"
$ws.Range("H2").Value = "双休,我在的时候不加班，后面不清楚有没有变化"
$ws.Range("I2").Value = "工资10%"
$ws.Range("J2").Value = "说是3个月，不确信可不可靠"
$ws.Range("K2").Value = "半年 不打折"
$ws.Range("L2").Value = "5k左右的thinkpad 可以自带，没补贴。"
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = "刷工牌打卡"
$ws.Range("O2").Value = ""
$ws.Range("P2").Value = ""
$ws.Range("Q2").Value = "2022-03-09 02:18:53"
$ws.Range("R2").Value = ""
$ws.Range("S2").Value = ""
